$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.834.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.675.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.838.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.69%  "

$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.51%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.780.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0515"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.89%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("E51").Value = "  +0.11%  "
